# BATCH_TEST_Batch.xlsx — "Lots of new changes"
# Append three new ShopText rows (TagTest entries) and reflect the resulting
# UI state (selection, column B auto-sized for the new long value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "test"
$ws.Range("B7").Value = "battle_text"
$ws.Range("C7").Value = "BattleText_TagTest"

$ws.Range("A8").Value = "test"
$ws.Range("B8").Value = "battle_text longest_name"
$ws.Range("C8").Value = "Multiple_TagTest"

$ws.Range("A9").Value = "test"
$ws.Range("B9").Value = "WWWWWWWWWWWWWWWWWWWWWWWWWWWWWWWWWWWW"
$ws.Range("C9").Value = "36Ws"

# Column B needs to widen to fit the new 36-character "W" string.
$ws.Columns.Item(2).ColumnWidth = 67

# Leave the selection where Excel would land after entering the data.
$ws.Range("C10").Select()

# Switch the calculation reference style to R1C1 (as in the authored edit).
$excel.ReferenceStyle = -4150
